$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sciwheel API note text was expanded with more detail about EndNote and timing
$ws.Range("E21").Value = "Sciwheel API works well. EndNote doesn't work. Took much longer to implement than 3 hours…"

# Time spent on this task increased from 1.25 to 3 hours
$ws.Range("D21").Value = 3

# Scroll the view down and move the selection to reflect where the user was working
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C21").Select()
